$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty date cells in rows 3 and 4 with the same
# "Fecha_Base" / "Fecha_Proyectada" values used in row 2, using the same
# date-formatted style (style index 2 -> numFmt 164, fillId 2) for all of
# them. This also makes the custom numFmtId=14 style (index 3, previously
# used by B4/C4) unused, so Excel drops it from cellXfs on save.
$ws.Range("B3").Value = 45992
$ws.Range("C3").Value = 46011
$ws.Range("B4").Value = 45992
$ws.Range("C4").Value = 46011

# B4/C4 previously used a different style (numFmtId 14) than B2/B3 (numFmtId 164).
# Make them consistent with the rest of the column by copying the style from B3/C3.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122) | Out-Null

# Restore the values after the paste-special formats-only operation (paste
# special with formats shouldn't touch values, but set them again just in case)
$ws.Range("B4").Value = 45992
$ws.Range("C4").Value = 46011

# Update the active selection to match the diff (C9 instead of D19)
$ws.Range("C9").Select()
